$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 411.33334
$ws.Range("I19").Value = 411.33334
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 411.33334
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -236.33334
$ws.Range("N19").ClearContents()
$ws.Range("H32").Value = 2465.611
$ws.Range("I32").Value = 2847.375
$ws.Range("J32").Value = 2160.2
$ws.Range("K32").Value = 2847.375
$ws.Range("L32").Value = 2160.2
$ws.Range("M32").Value = -2521.375
$ws.Range("N32").Value = -2812.2
$ws.Range("H41").Value = 279.55554
$ws.Range("I41").Value = 219.71428
$ws.Range("J41").Value = 317.63635
$ws.Range("K41").Value = 219.71428
$ws.Range("L41").Value = 317.63635
$ws.Range("M41").Value = 220.28572
$ws.Range("N41").Value = -1197.63635
$ws.Range("H51").Value = 3650.3333
$ws.Range("I51").Value = 2200
$ws.Range("J51").Value = 3940.4
$ws.Range("K51").Value = 2200
$ws.Range("L51").Value = 3940.4
$ws.Range("M51").Value = -1716
$ws.Range("N51").Value = -4908.4
$ws.Range("H62").Value = 2167.5
$ws.Range("I62").Value = 1801
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 1801
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -1177
$ws.Range("N62").Value = -5248
$ws.Range("H65").Value = 2167.5
$ws.Range("I65").Value = 1801
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 9005
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -5885
$ws.Range("N65").Value = -26240
$ws.Range("H98").Value = 3405.875
$ws.Range("I98").Value = 1608.5454
$ws.Range("J98").Value = 7360
$ws.Range("K98").Value = 1608.5454
$ws.Range("L98").Value = 7360
$ws.Range("M98").Value = -110.5454
$ws.Range("N98").Value = -10356
$ws.Range("H116").Value = 1948.0646
$ws.Range("I116").Value = 1446.6666
$ws.Range("K116").Value = 1446.6666
$ws.Range("M116").Value = 1995.3334
$ws.Range("H122").Value = 3405.875
$ws.Range("I122").Value = 1608.5454
$ws.Range("J122").Value = 7360
$ws.Range("K122").Value = 4825.6362
$ws.Range("L122").Value = 22080
$ws.Range("M122").Value = -2375.6362
$ws.Range("N122").Value = -26980
$ws.Range("H132").Value = 1957358.5
$ws.Range("I132").Value = 2408414
$ws.Range("J132").Value = 2783.9167
$ws.Range("K132").Value = 7225242
$ws.Range("L132").Value = 8351.750100000001
$ws.Range("M132").Value = -7222712
$ws.Range("N132").Value = -13411.7501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 11182
$ws.Range("I61").Value = 11182
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 11182
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -10970
$ws.Range("N61").ClearContents()
$ws.Range("H132").Value = 35470.832
$ws.Range("I132").Value = 35270.668
$ws.Range("J132").Value = 35671
$ws.Range("K132").Value = 105812.004
$ws.Range("L132").Value = 107013
$ws.Range("M132").Value = -103282.004
$ws.Range("N132").Value = -112073
$ws.Range("H136").Value = 11182
$ws.Range("I136").Value = 11182
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 33546
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -30996
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1253.5
$ws.Range("I107").Value = 1130.25
$ws.Range("J107").Value = 1500
$ws.Range("K107").Value = 1130.25
$ws.Range("L107").Value = 1500
$ws.Range("M107").Value = 789.75
$ws.Range("N107").Value = -5340

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 5421.515
$ws.Range("I132").Value = 6316.1304
$ws.Range("J132").Value = 3363.9
$ws.Range("K132").Value = 18948.3912
$ws.Range("L132").Value = 10091.7
$ws.Range("M132").Value = -16418.3912
$ws.Range("N132").Value = -15151.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2862.5
$ws.Range("I80").Value = 2816.6667
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2816.6667
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1818.6667
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2862.5
$ws.Range("I83").Value = 2816.6667
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 14083.3335
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -9091.333500000001
$ws.Range("N83").Value = -24984
$ws.Range("H97").Value = 1503.6666
$ws.Range("I97").Value = 2000
$ws.Range("J97").Value = 1255.5
$ws.Range("K97").Value = 2000
$ws.Range("L97").Value = 1255.5
$ws.Range("M97").Value = -1504
$ws.Range("N97").Value = -2247.5
$ws.Range("H107").Value = 2024737.4
$ws.Range("I107").Value = 155.3
$ws.Range("J107").Value = 3290101.2
$ws.Range("K107").Value = 155.3
$ws.Range("L107").Value = 3290101.2
$ws.Range("M107").Value = 1764.7
$ws.Range("N107").Value = -3293941.2
$ws.Range("H126").Value = 1618.75
$ws.Range("I126").Value = 1658.8572
$ws.Range("J126").Value = 1587.5555
$ws.Range("K126").Value = 4976.571599999999
$ws.Range("L126").Value = 4762.666499999999
$ws.Range("M126").Value = -2506.571599999999
$ws.Range("N126").Value = -9702.666499999999
$ws.Range("H132").Value = 17661.385
$ws.Range("I132").Value = 14125
$ws.Range("J132").Value = 23319.6
$ws.Range("K132").Value = 42375
$ws.Range("L132").Value = 69958.79999999999
$ws.Range("M132").Value = -39845
$ws.Range("N132").Value = -75018.79999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2182.1428
$ws.Range("I61").Value = 1866.1428
$ws.Range("J61").Value = 2814.1428
$ws.Range("K61").Value = 1866.1428
$ws.Range("L61").Value = 2814.1428
$ws.Range("M61").Value = -1664.1428
$ws.Range("N61").Value = -3218.1428
$ws.Range("H82").Value = 2799.9333
$ws.Range("I82").Value = 1899.5
$ws.Range("J82").Value = 2938.4614
$ws.Range("K82").Value = 1899.5
$ws.Range("L82").Value = 2938.4614
$ws.Range("M82").Value = -1538.5
$ws.Range("N82").Value = -3660.4614
$ws.Range("H85").Value = 2799.9333
$ws.Range("I85").Value = 1899.5
$ws.Range("J85").Value = 2938.4614
$ws.Range("K85").Value = 1899.5
$ws.Range("L85").Value = 2938.4614
$ws.Range("M85").Value = -651.5
$ws.Range("N85").Value = -5434.4614
$ws.Range("H113").Value = 2182.1428
$ws.Range("I113").Value = 1866.1428
$ws.Range("J113").Value = 2814.1428
$ws.Range("K113").Value = 1866.1428
$ws.Range("L113").Value = 2814.1428
$ws.Range("M113").Value = 303.8571999999999
$ws.Range("N113").Value = -7154.1428
$ws.Range("H132").Value = 94336.17999999999
$ws.Range("I132").Value = 126961.875
$ws.Range("K132").Value = 380885.625
$ws.Range("M132").Value = -378355.625
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H136").Value = 3080.7058
$ws.Range("I136").Value = 2417.842
$ws.Range("J136").Value = 3920.3333
$ws.Range("K136").Value = 7253.526
$ws.Range("L136").Value = 11760.9999
$ws.Range("M136").Value = -4703.526
$ws.Range("N136").Value = -16860.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5538.5312
$ws.Range("I132").Value = 5863.0435
$ws.Range("J132").Value = 4709.222
$ws.Range("K132").Value = 17589.1305
$ws.Range("L132").Value = 14127.666
$ws.Range("M132").Value = -15059.1305
$ws.Range("N132").Value = -19187.666
$ws.Range("H136").Value = 34487550
$ws.Range("I136").Value = 38466410
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 115399230
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -115396680
$ws.Range("N136").Value = -17599.9995
